# Auto-generated Excel COM-interop script to apply financial model updates
# to the Goblin Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H15").Value = 507.73468
$ws.Range("I15").Value = 507.73468
$ws.Range("K15").Value = 1523.20404
$ws.Range("M15").Value = -1354.20404
$ws.Range("H19").Value = 656.8570999999999
$ws.Range("I19").Value = 874.75
$ws.Range("J19").Value = 366.33334
$ws.Range("K19").Value = 874.75
$ws.Range("L19").Value = 366.33334
$ws.Range("M19").Value = -699.75
$ws.Range("N19").Value = -716.33334
$ws.Range("H39").Value = 421.84616
$ws.Range("J39").Value = 1004.6667
$ws.Range("L39").Value = 3014.0001
$ws.Range("N39").Value = -3606.0001
$ws.Range("H40").Value = 2318.182
$ws.Range("J40").Value = 3066.5
$ws.Range("L40").Value = 3066.5
$ws.Range("N40").Value = -3416.5
$ws.Range("H41").Value = 1221.4286
$ws.Range("I41").Value = 915.4
$ws.Range("J41").Value = 1986.5
$ws.Range("K41").Value = 915.4
$ws.Range("L41").Value = 1986.5
$ws.Range("M41").Value = -475.4
$ws.Range("N41").Value = -2866.5
$ws.Range("H46").Value = 45125
$ws.Range("I46").Value = 60000
$ws.Range("J46").Value = 500
$ws.Range("K46").Value = 180000
$ws.Range("L46").Value = 1500
$ws.Range("M46").Value = -179881
$ws.Range("N46").Value = -1738
$ws.Range("H60").Value = 45125
$ws.Range("I60").Value = 60000
$ws.Range("J60").Value = 500
$ws.Range("K60").Value = 180000
$ws.Range("L60").Value = 1500
$ws.Range("M60").Value = -179516
$ws.Range("N60").Value = -2468
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H86").Value = 3908.25
$ws.Range("I86").Value = 3999.4443
$ws.Range("K86").Value = 3999.4443
$ws.Range("M86").Value = -2876.4443
$ws.Range("H87").Value = 74999.2
$ws.Range("J87").Value = 74999.2
$ws.Range("L87").Value = 74999.2
$ws.Range("N87").Value = -77495.2
$ws.Range("H88").Value = 6293.4
$ws.Range("I88").Value = 3467
$ws.Range("J88").Value = 7000
$ws.Range("K88").Value = 3467
$ws.Range("L88").Value = 7000
$ws.Range("M88").Value = -3061
$ws.Range("N88").Value = -7812
$ws.Range("H89").Value = 3908.25
$ws.Range("I89").Value = 3999.4443
$ws.Range("K89").Value = 19997.2215
$ws.Range("M89").Value = -14381.2215
$ws.Range("H90").Value = 74999.2
$ws.Range("J90").Value = 74999.2
$ws.Range("L90").Value = 224997.6
$ws.Range("N90").Value = -237477.6
$ws.Range("H91").Value = 6293.4
$ws.Range("I91").Value = 3467
$ws.Range("J91").Value = 7000
$ws.Range("K91").Value = 3467
$ws.Range("L91").Value = 7000
$ws.Range("M91").Value = -2063
$ws.Range("N91").Value = -9808
$ws.Range("H101").Value = 1411.8462
$ws.Range("J101").Value = 15000
$ws.Range("L101").Value = 45000
$ws.Range("N101").Value = -48244
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H103").Value = 1068.9231
$ws.Range("J103").Value = 1144.3334
$ws.Range("L103").Value = 3433.0002
$ws.Range("N103").Value = -4605.0002
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H132").Value = 2875.9333
$ws.Range("I132").Value = 2366.6428
$ws.Range("K132").Value = 7099.928400000001
$ws.Range("M132").Value = -4569.928400000001
$ws.Range("H138").Value = 4096.7896
$ws.Range("J138").Value = 4864
$ws.Range("L138").Value = 14592
$ws.Range("N138").Value = -24872

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2433
$ws.Range("I2").Value = 854.5714
$ws.Range("K2").Value = 854.5714
$ws.Range("M2").Value = -741.5714
$ws.Range("H32").Value = 1850.9028
$ws.Range("I32").Value = 1771.7537
$ws.Range("K32").Value = 1771.7537
$ws.Range("M32").Value = -1484.7537
$ws.Range("H61").Value = 3680.0789
$ws.Range("I61").Value = 3592
$ws.Range("K61").Value = 3592
$ws.Range("M61").Value = -3380
$ws.Range("H74").Value = 3548.0557
$ws.Range("J74").Value = 3445.25
$ws.Range("L74").Value = 3445.25
$ws.Range("N74").Value = -5193.25
$ws.Range("H77").Value = 3548.0557
$ws.Range("J77").Value = 3445.25
$ws.Range("L77").Value = 17226.25
$ws.Range("N77").Value = -25962.25
$ws.Range("H102").Value = 3932.9092
$ws.Range("I102").Value = 2207.2942
$ws.Range("K102").Value = 2207.2942
$ws.Range("M102").Value = -585.2941999999998
$ws.Range("H116").Value = 2433
$ws.Range("I116").Value = 854.5714
$ws.Range("K116").Value = 854.5714
$ws.Range("M116").Value = 1439.4286
$ws.Range("H132").Value = 3468.6086
$ws.Range("I132").Value = 3421.7727
$ws.Range("K132").Value = 10265.3181
$ws.Range("M132").Value = -7735.3181
$ws.Range("H136").Value = 3680.0789
$ws.Range("I136").Value = 3592
$ws.Range("K136").Value = 10776
$ws.Range("M136").Value = -8226

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2433
$ws.Range("I3").Value = 854.5714
$ws.Range("K3").Value = 854.5714
$ws.Range("M3").Value = -740.5714
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H135").Value = 78000
$ws.Range("J135").Value = 78000
$ws.Range("L135").Value = 78000
$ws.Range("N135").Value = -88140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2660.8572
$ws.Range("I134").Value = 2836.4167
$ws.Range("K134").Value = 8509.250100000001
$ws.Range("M134").Value = -5974.250100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 4629.4165
$ws.Range("J137").Value = 5566.4
$ws.Range("L137").Value = 16699.2
$ws.Range("N137").Value = -26899.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 4190.8438
$ws.Range("I97").Value = 1246.9231
$ws.Range("J97").Value = 16947.834
$ws.Range("K97").Value = 1246.9231
$ws.Range("L97").Value = 16947.834
$ws.Range("M97").Value = -750.9231
$ws.Range("N97").Value = -17939.834

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 50650
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H46").Value = 3444.3
$ws.Range("I46").Value = 2991.6667
$ws.Range("J46").Value = 4123.25
$ws.Range("K46").Value = 2991.6667
$ws.Range("L46").Value = 4123.25
$ws.Range("M46").Value = -2803.6667
$ws.Range("N46").Value = -4499.25
$ws.Range("H61").Value = 4677
$ws.Range("I61").Value = 2891.6924
$ws.Range("J61").Value = 6334.7856
$ws.Range("K61").Value = 2891.6924
$ws.Range("L61").Value = 6334.7856
$ws.Range("M61").Value = -2689.6924
$ws.Range("N61").Value = -6738.7856
$ws.Range("H68").Value = 5168.2856
$ws.Range("I68").Value = 2249.7778
$ws.Range("K68").Value = 2249.7778
$ws.Range("M68").Value = -1500.7778
$ws.Range("H71").Value = 5168.2856
$ws.Range("I71").Value = 2249.7778
$ws.Range("K71").Value = 11248.889
$ws.Range("M71").Value = -7504.888999999999
$ws.Range("H113").Value = 4677
$ws.Range("I113").Value = 2891.6924
$ws.Range("J113").Value = 6334.7856
$ws.Range("K113").Value = 2891.6924
$ws.Range("L113").Value = 6334.7856
$ws.Range("M113").Value = -721.6923999999999
$ws.Range("N113").Value = -10674.7856
$ws.Range("H138").Value = 83330.664
$ws.Range("J138").Value = 83330.664
$ws.Range("L138").Value = 83330.664
$ws.Range("N138").Value = -93610.664
$ws.Range("H140").Value = 40429
$ws.Range("J140").Value = 40429
$ws.Range("L140").Value = 40429
$ws.Range("N140").Value = -50789

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 200000
$ws.Range("I52").Value = 200000
$ws.Range("K52").Value = 200000
$ws.Range("M52").Value = -199774
$ws.Range("H81").Value = 2500.3333
$ws.Range("J81").Value = 5000
$ws.Range("L81").Value = 10000
$ws.Range("N81").Value = -12122
$ws.Range("H84").Value = 2500.3333
$ws.Range("J84").Value = 5000
$ws.Range("L84").Value = 50000
$ws.Range("N84").Value = -60608
$ws.Range("H113").Value = 946.9459000000001
$ws.Range("I113").Value = 834.4643
$ws.Range("J113").Value = 1296.8889
$ws.Range("K113").Value = 2503.3929
$ws.Range("L113").Value = 3890.6667
$ws.Range("M113").Value = -333.3928999999998
$ws.Range("N113").Value = -8230.6667
$ws.Range("H122").Value = 5080.269
$ws.Range("J122").Value = 9033.111000000001
$ws.Range("L122").Value = 27099.333
$ws.Range("N122").Value = -31999.333
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H126").Value = 3309.7693
$ws.Range("I126").Value = 3309.7693
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9929.3079
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -7459.3079
$ws.Range("N126").ClearContents()
